$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.371.10'
$ws.Range("E2").Value = '  -0.20%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.522.15'
$ws.Range("E3").Value = '  +1.26%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '521.74'
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.33'
$ws.Range("E6").Value = '  -1.64%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.520.20'
$ws.Range("E9").Value = '  +0.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0974'
$ws.Range("E10").Value = '  -1.73%  '
$ws.Range("E11").Value = '  -1.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.16'
$ws.Range("E12").Value = '  -3.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.331'
$ws.Range("E13").Value = '  -2.90%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.970.28'
$ws.Range("E14").Value = '  +1.27%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.350.11'
$ws.Range("E15").Value = '  -1.51%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.09'
$ws.Range("E16").Value = '  -0.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000134'
$ws.Range("E17").Value = '  -1.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.531.10'
$ws.Range("E18").Value = '  +1.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.65'
$ws.Range("E19").Value = '  -0.96%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '321.67'
$ws.Range("E20").Value = '  -0.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.15'
$ws.Range("E21").Value = '  -1.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.18'
$ws.Range("E22").Value = '  +6.85%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.63'
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("E25").Value = '  -1.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.159'
$ws.Range("E27").Value = '  -1.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.39'
$ws.Range("E28").Value = '  -0.62%  '
$ws.Range("E29").Value = '  -0.63%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '168.46'
$ws.Range("E30").Value = '  -0.78%  '
$ws.Range("E31").Value = '  +0.66%  '
$ws.Range("E32").Value = '  -0.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.28'
$ws.Range("E33").Value = '  -1.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.17'
$ws.Range("E36").Value = '  -0.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.25'
$ws.Range("E37").Value = '  -7.04%  '
$ws.Range("E38").Value = '  -3.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.49'
$ws.Range("E39").Value = '  +0.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.40'
$ws.Range("E40").Value = '  -0.84%  '
$ws.Range("E41").Value = '  -4.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '276.25'
$ws.Range("E42").Value = '  -1.17%  '
$ws.Range("E43").Value = '  -1.24%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '130.45'
$ws.Range("E44").Value = '  +4.50%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.99'
$ws.Range("E45").Value = '  -4.55%  '
$ws.Range("E46").Value = '  -0.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0918'
$ws.Range("E47").Value = '  +0.63%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0499'
$ws.Range("E48").Value = '  +0.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '17.72'
$ws.Range("E49").Value = '  -0.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0214'
$ws.Range("E50").Value = '  -0.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.90'
$ws.Range("E51").Value = '  -1.81%  '
